$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("district_price_per_m2")
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = 'Reinickendorf'
$ws.Cells.Item(2, 3).Value = 31.16170698273666
$ws.Cells.Item(2, 4).Value = 30.625
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Lichtenberg'
$ws.Cells.Item(3, 3).Value = 30.3144676245913
$ws.Cells.Item(3, 4).Value = 30
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = 'Steglitz-Zehlendorf'
$ws.Cells.Item(4, 3).Value = 30.07451461223216
$ws.Cells.Item(4, 4).Value = 30
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'Marzahn-Hellersdorf'
$ws.Cells.Item(5, 3).Value = 29.48557408791845
$ws.Cells.Item(5, 4).Value = 29.41176470588235
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'Mitte'
$ws.Cells.Item(6, 3).Value = 31.08476283926247
$ws.Cells.Item(6, 4).Value = 29.16666666666667
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = 'Spandau'
$ws.Cells.Item(7, 3).Value = 30.6194903245398
$ws.Cells.Item(7, 4).Value = 29.16666666666667
$ws.Cells.Item(8, 1).Value = 12
$ws.Cells.Item(8, 2).Value = 'Berlin Insgesamt'
$ws.Cells.Item(8, 3).Value = 29.9776885198542
$ws.Cells.Item(8, 4).Value = 29.02850877192982
$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = 'Pankow'
$ws.Cells.Item(9, 3).Value = 29.92796413217931
$ws.Cells.Item(9, 4).Value = 28.89035087719298
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = 'Treptow-Köpenick'
$ws.Cells.Item(10, 3).Value = 30.68054449536447
$ws.Cells.Item(10, 4).Value = 28.84615384615385
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Tempelhof-Schöneberg'
$ws.Cells.Item(11, 3).Value = 28.93973636828191
$ws.Cells.Item(11, 4).Value = 28.29710144927536
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = 'Friedrichshain-Kreuzberg'
$ws.Cells.Item(12, 3).Value = 29.35592648441228
$ws.Cells.Item(12, 4).Value = 28.09523809523809
$ws.Cells.Item(13, 1).Value = 0
$ws.Cells.Item(13, 2).Value = 'Charlottenburg-Wilmersdorf'
$ws.Cells.Item(13, 3).Value = 30.66872772425488
$ws.Cells.Item(13, 4).Value = 28.04761904761905
$ws.Cells.Item(14, 1).Value = 5
$ws.Cells.Item(14, 2).Value = 'Neukölln'
$ws.Cells.Item(14, 3).Value = 27.41884656247671
$ws.Cells.Item(14, 4).Value = 26.66666666666667

$ws = $wb.Worksheets.Item("district_adjusted_price_per_m2")
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = 'Reinickendorf'
$ws.Cells.Item(2, 3).Value = 19.67619550281314
$ws.Cells.Item(2, 4).Value = 20
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = 'Spandau'
$ws.Cells.Item(3, 3).Value = 19.03477995695436
$ws.Cells.Item(3, 4).Value = 19.51364479202366
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = 'Steglitz-Zehlendorf'
$ws.Cells.Item(4, 3).Value = 18.598314975219
$ws.Cells.Item(4, 4).Value = 18.67518375761195
$ws.Cells.Item(5, 1).Value = 0
$ws.Cells.Item(5, 2).Value = 'Charlottenburg-Wilmersdorf'
$ws.Cells.Item(5, 3).Value = 19.13714138598957
$ws.Cells.Item(5, 4).Value = 18.33700114681742
$ws.Cells.Item(6, 1).Value = 6
$ws.Cells.Item(6, 2).Value = 'Pankow'
$ws.Cells.Item(6, 3).Value = 18.82573745548559
$ws.Cells.Item(6, 4).Value = 18.15456181016873
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = 'Mitte'
$ws.Cells.Item(7, 3).Value = 18.95857842330486
$ws.Cells.Item(7, 4).Value = 17.85243145838252
$ws.Cells.Item(8, 1).Value = 12
$ws.Cells.Item(8, 2).Value = 'Berlin Insgesamt'
$ws.Cells.Item(8, 3).Value = 18.42986563272193
$ws.Cells.Item(8, 4).Value = 17.79718347112674
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = 'Lichtenberg'
$ws.Cells.Item(9, 3).Value = 18.46710207296916
$ws.Cells.Item(9, 4).Value = 17.74193548387097
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = 'Tempelhof-Schöneberg'
$ws.Cells.Item(10, 3).Value = 18.3573921852684
$ws.Cells.Item(10, 4).Value = 17.47003829523662
$ws.Cells.Item(11, 1).Value = 3
$ws.Cells.Item(11, 2).Value = 'Marzahn-Hellersdorf'
$ws.Cells.Item(11, 3).Value = 17.66405030947761
$ws.Cells.Item(11, 4).Value = 17.25384615384615
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Treptow-Köpenick'
$ws.Cells.Item(12, 3).Value = 17.88374849583476
$ws.Cells.Item(12, 4).Value = 17.22833925957488
$ws.Cells.Item(13, 1).Value = 5
$ws.Cells.Item(13, 2).Value = 'Neukölln'
$ws.Cells.Item(13, 3).Value = 17.15581903895128
$ws.Cells.Item(13, 4).Value = 16.73453983057948
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = 'Friedrichshain-Kreuzberg'
$ws.Cells.Item(14, 3).Value = 17.39952779039545
$ws.Cells.Item(14, 4).Value = 16.61828875316738

$ws = $wb.Worksheets.Item("district_price_per_m2_extremes")
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = 'Marzahn-Hellersdorf'
$ws.Cells.Item(2, 3).Value = 0.6666666666666666
$ws.Cells.Item(2, 4).Value = 57.91666666666666
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = 'Neukölln'
$ws.Cells.Item(3, 3).Value = 0.7
$ws.Cells.Item(3, 4).Value = 81.14285714285714
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 'Friedrichshain-Kreuzberg'
$ws.Cells.Item(4, 3).Value = 0.75
$ws.Cells.Item(4, 4).Value = 85
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Mitte'
$ws.Cells.Item(5, 3).Value = 0.8571428571428571
$ws.Cells.Item(5, 4).Value = 110
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = 'Tempelhof-Schöneberg'
$ws.Cells.Item(6, 3).Value = 0.875
$ws.Cells.Item(6, 4).Value = 63.33333333333334
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Pankow'
$ws.Cells.Item(7, 3).Value = 0.9090909090909091
$ws.Cells.Item(7, 4).Value = 100
$ws.Cells.Item(8, 1).Value = 11
$ws.Cells.Item(8, 2).Value = 'Treptow-Köpenick'
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 74.2
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Spandau'
$ws.Cells.Item(9, 3).Value = 1.428571428571429
$ws.Cells.Item(9, 4).Value = 56.81818181818182
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = 'Reinickendorf'
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 70
$ws.Cells.Item(11, 1).Value = 0
$ws.Cells.Item(11, 2).Value = 'Charlottenburg-Wilmersdorf'
$ws.Cells.Item(11, 3).Value = 2.125
$ws.Cells.Item(11, 4).Value = 105
$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = 'Lichtenberg'
$ws.Cells.Item(12, 3).Value = 2.142857142857143
$ws.Cells.Item(12, 4).Value = 84.0909090909091
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = 'Steglitz-Zehlendorf'
$ws.Cells.Item(13, 3).Value = 2.166666666666667
$ws.Cells.Item(13, 4).Value = 57

$ws = $wb.Worksheets.Item("district_adjusted_price_per_m2_extremes")
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = 'Neukölln'
$ws.Cells.Item(2, 3).Value = 0.5117904641168637
$ws.Cells.Item(2, 4).Value = 43.6875
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = 'Mitte'
$ws.Cells.Item(3, 3).Value = 0.5305124069678346
$ws.Cells.Item(3, 4).Value = 62.85714285714285
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = 'Pankow'
$ws.Cells.Item(4, 3).Value = 0.5821908875751723
$ws.Cells.Item(4, 4).Value = 61.81818181818182
$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = 'Spandau'
$ws.Cells.Item(5, 3).Value = 0.5868804139798212
$ws.Cells.Item(5, 4).Value = 41.1764705882353
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 'Friedrichshain-Kreuzberg'
$ws.Cells.Item(6, 3).Value = 0.715559793334368
$ws.Cells.Item(6, 4).Value = 46
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = 'Treptow-Köpenick'
$ws.Cells.Item(7, 3).Value = 0.7368421052631579
$ws.Cells.Item(7, 4).Value = 65
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = 'Tempelhof-Schöneberg'
$ws.Cells.Item(8, 3).Value = 0.8542691130192367
$ws.Cells.Item(8, 4).Value = 59.375
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'Reinickendorf'
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 45
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = 'Lichtenberg'
$ws.Cells.Item(10, 3).Value = 1.025641025641026
$ws.Cells.Item(10, 4).Value = 44.98587593420606
$ws.Cells.Item(11, 1).Value = 0
$ws.Cells.Item(11, 2).Value = 'Charlottenburg-Wilmersdorf'
$ws.Cells.Item(11, 3).Value = 1.666666666666667
$ws.Cells.Item(11, 4).Value = 42.35294117647059
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = 'Steglitz-Zehlendorf'
$ws.Cells.Item(12, 3).Value = 2.329378839231888
$ws.Cells.Item(12, 4).Value = 37.22222222222222
$ws.Cells.Item(13, 1).Value = 3
$ws.Cells.Item(13, 2).Value = 'Marzahn-Hellersdorf'
$ws.Cells.Item(13, 3).Value = 5.866666666666666
$ws.Cells.Item(13, 4).Value = 34.75

$ws = $wb.Worksheets.Item("Price Comparison")
$ws.Cells.Item(2, 1).Value = 'Spandau'
$ws.Cells.Item(2, 2).Value = 29.16666666666667
$ws.Cells.Item(2, 3).Value = 19.51364479202366
$ws.Cells.Item(2, 4).Value = 9.341384152837394
$ws.Cells.Item(2, 5).Value = 19.82528251382927
$ws.Cells.Item(2, 6).Value = 10.17226063918627
$ws.Cells.Item(2, 7).Value = 312.230673628891
$ws.Cells.Item(2, 8).Value = 208.8945757154897
$ws.Cells.Item(3, 1).Value = 'Reinickendorf'
$ws.Cells.Item(3, 2).Value = 30.625
$ws.Cells.Item(3, 3).Value = 20
$ws.Cells.Item(3, 4).Value = 10.47216271058341
$ws.Cells.Item(3, 5).Value = 20.1528372894166
$ws.Cells.Item(3, 6).Value = 9.527837289416594
$ws.Cells.Item(3, 7).Value = 292.4419801943077
$ws.Cells.Item(3, 8).Value = 190.9825176779152
$ws.Cells.Item(4, 1).Value = 'Marzahn-Hellersdorf'
$ws.Cells.Item(4, 2).Value = 29.41176470588235
$ws.Cells.Item(4, 3).Value = 17.25384615384615
$ws.Cells.Item(4, 4).Value = 10.3678218974249
$ws.Cells.Item(4, 5).Value = 19.04394280845745
$ws.Cells.Item(4, 6).Value = 6.886024256421255
$ws.Cells.Item(4, 7).Value = 283.6831592678832
$ws.Cells.Item(4, 8).Value = 166.4172699391332
$ws.Cells.Item(5, 1).Value = 'Lichtenberg'
$ws.Cells.Item(5, 2).Value = 30
$ws.Cells.Item(5, 3).Value = 17.74193548387097
$ws.Cells.Item(5, 4).Value = 11.1179280113382
$ws.Cells.Item(5, 5).Value = 18.8820719886618
$ws.Cells.Item(5, 6).Value = 6.624007472532766
$ws.Cells.Item(5, 7).Value = 269.8344508923391
$ws.Cells.Item(5, 8).Value = 159.5795139685877
$ws.Cells.Item(6, 1).Value = 'Neukölln'
$ws.Cells.Item(6, 2).Value = 26.66666666666667
$ws.Cells.Item(6, 3).Value = 16.73453983057948
$ws.Cells.Item(6, 4).Value = 11.52667883982897
$ws.Cells.Item(6, 5).Value = 15.1399878268377
$ws.Cells.Item(6, 6).Value = 5.207860990750508
$ws.Cells.Item(6, 7).Value = 231.3473554457282
$ws.Cells.Item(6, 8).Value = 145.1809325402162
$ws.Cells.Item(7, 1).Value = 'Tempelhof-Schöneberg'
$ws.Cells.Item(7, 2).Value = 28.29710144927536
$ws.Cells.Item(7, 3).Value = 17.47003829523662
$ws.Cells.Item(7, 4).Value = 12.13993599794491
$ws.Cells.Item(7, 5).Value = 16.15716545133045
$ws.Cells.Item(7, 6).Value = 5.330102297291713
$ws.Cells.Item(7, 7).Value = 233.0910266253923
$ws.Cells.Item(7, 8).Value = 143.9055222218141
$ws.Cells.Item(8, 1).Value = 'Berlin Insgesamt'
$ws.Cells.Item(8, 2).Value = 29.02850877192982
$ws.Cells.Item(8, 3).Value = 17.79718347112674
$ws.Cells.Item(8, 4).Value = 12.50626454510267
$ws.Cells.Item(8, 5).Value = 16.52224422682715
$ws.Cells.Item(8, 6).Value = 5.290918926024078
$ws.Cells.Item(8, 7).Value = 232.1117442161985
$ws.Cells.Item(8, 8).Value = 142.3061491058571
$ws.Cells.Item(9, 1).Value = 'Steglitz-Zehlendorf'
$ws.Cells.Item(9, 2).Value = 30
$ws.Cells.Item(9, 3).Value = 18.67518375761195
$ws.Cells.Item(9, 4).Value = 13.23848794563714
$ws.Cells.Item(9, 5).Value = 16.76151205436286
$ws.Cells.Item(9, 6).Value = 5.436695811974817
$ws.Cells.Item(9, 7).Value = 226.6119826009795
$ws.Cells.Item(9, 8).Value = 141.0673472250018
$ws.Cells.Item(10, 1).Value = 'Pankow'
$ws.Cells.Item(10, 2).Value = 28.89035087719298
$ws.Cells.Item(10, 3).Value = 18.15456181016873
$ws.Cells.Item(10, 4).Value = 13.63586323304723
$ws.Cells.Item(10, 5).Value = 15.25448764414575
$ws.Cells.Item(10, 6).Value = 4.518698577121501
$ws.Cells.Item(10, 7).Value = 211.8703479452309
$ws.Cells.Item(10, 8).Value = 133.1383389514365
$ws.Cells.Item(11, 1).Value = 'Treptow-Köpenick'
$ws.Cells.Item(11, 2).Value = 28.84615384615385
$ws.Cells.Item(11, 3).Value = 17.22833925957488
$ws.Cells.Item(11, 4).Value = 13.02911475531228
$ws.Cells.Item(11, 5).Value = 15.81703909084156
$ws.Cells.Item(11, 6).Value = 4.199224504262597
$ws.Cells.Item(11, 7).Value = 221.3976497090301
$ws.Cells.Item(11, 8).Value = 132.2295457759359
$ws.Cells.Item(12, 1).Value = 'Charlottenburg-Wilmersdorf'
$ws.Cells.Item(12, 2).Value = 28.04761904761905
$ws.Cells.Item(12, 3).Value = 18.33700114681742
$ws.Cells.Item(12, 4).Value = 16.42457502011024
$ws.Cells.Item(12, 5).Value = 11.6230440275088
$ws.Cells.Item(12, 6).Value = 1.912426126707178
$ws.Cells.Item(12, 7).Value = 170.7661782011258
$ws.Cells.Item(12, 8).Value = 111.6436871235061
$ws.Cells.Item(13, 1).Value = 'Mitte'
$ws.Cells.Item(13, 2).Value = 29.16666666666667
$ws.Cells.Item(13, 3).Value = 17.85243145838252
$ws.Cells.Item(13, 4).Value = 17.57533810105219
$ws.Cells.Item(13, 5).Value = 11.59132856561448
$ws.Cells.Item(13, 6).Value = 0.2770933573303331
$ws.Cells.Item(13, 7).Value = 165.9522365883848
$ws.Cells.Item(13, 8).Value = 101.5766032820373
$ws.Cells.Item(14, 1).Value = 'Friedrichshain-Kreuzberg'
$ws.Cells.Item(14, 2).Value = 28.09523809523809
$ws.Cells.Item(14, 3).Value = 16.61828875316738
$ws.Cells.Item(14, 4).Value = 16.43813941562651
$ws.Cells.Item(14, 5).Value = 11.65709867961159
$ws.Cells.Item(14, 6).Value = 0.1801493375408754
$ws.Cells.Item(14, 7).Value = 170.9149520202393
$ws.Cells.Item(14, 8).Value = 101.0959229203861

